$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("C3").Value = -11.4253
$ws.Range("C4").Value = -12.34770000000001
$ws.Range("B7").Value = 5.579099999999998
$ws.Range("A8").Value = -22.47590000000001
$ws.Range("A10").Value = -21.2311
$ws.Range("D10").Value = -7.856899999999999
$ws.Range("C11").Value = -11.6147
$ws.Range("A12").Value = -21.59010000000001
$ws.Range("D12").Value = -7.317299999999989
$ws.Range("D13").Value = -9.224599999999995
$ws.Range("B14").Value = 5.888000000000002
$ws.Range("C14").Value = -12.29589999999999
$ws.Range("D14").Value = -7.658000000000004
$ws.Range("B15").Value = 4.796399999999995
$ws.Range("E16").Value = 16.50200000000001
$ws.Range("A18").Value = -21.82469999999999
$ws.Range("B18").Value = 5.981599999999995
$ws.Range("C18").Value = -11.4896
$ws.Range("E18").Value = 18.36460000000003
$ws.Range("C19").Value = -11.91390000000001
$ws.Range("B20").Value = 8.660200000000001
$ws.Range("C21").Value = -12.805
$ws.Range("E21").Value = 16.68700000000002
$ws.Range("E22").Value = 17.2312
$ws.Range("A25").Value = -21.4696
$ws.Range("E26").Value = 16.19799999999999
$ws.Range("C27").Value = -12.28689999999999
$ws.Range("E27").Value = 16.81659999999999
$ws.Range("B29").Value = 4.920200000000003
$ws.Range("D29").Value = -7.430899999999993
$ws.Range("B30").Value = 5.525400000000003
$ws.Range("B31").Value = 4.124499999999997
$ws.Range("C31").Value = -13.2457
$ws.Range("D32").Value = -8.773799999999994
$ws.Range("B35").Value = 8.241300000000008
$ws.Range("D35").Value = -8.394999999999998
$ws.Range("A37").Value = -19.3021
$ws.Range("C38").Value = -12.67959999999999
$ws.Range("E39").Value = 16.2189
$ws.Range("B40").Value = 9.226199999999997
$ws.Range("C42").Value = -12.12380000000001
$ws.Range("D43").Value = -8.126799999999998
$ws.Range("B44").Value = 5.110900000000004
$ws.Range("C44").Value = -12.92859999999999
$ws.Range("E44").Value = 16.6186
$ws.Range("C47").Value = -12.2826
$ws.Range("D48").Value = -7.6934
$ws.Range("D49").Value = -8.407400000000003
$ws.Range("B50").Value = 5.406499999999998
$ws.Range("D50").Value = -8.067699999999993
$ws.Range("D51").Value = -8.606599999999997
$ws.Range("E51").Value = 16.5408
$ws.Range("B54").Value = 4.7938
$ws.Range("E54").Value = 16.55000000000001
$ws.Range("A55").Value = -22.56080000000001
$ws.Range("C56").Value = -13.17379999999999
$ws.Range("D56").Value = -8.295000000000002
$ws.Range("E57").Value = 16.80599999999999
$ws.Range("C58").Value = -13.5976
$ws.Range("E58").Value = 16.51040000000001
$ws.Range("E60").Value = 15.80470000000002
$ws.Range("D61").Value = -8.043499999999996
$ws.Range("E63").Value = 17.72920000000001
$ws.Range("C65").Value = -12.5445
$ws.Range("A68").Value = -21.48170000000001
$ws.Range("B68").Value = 4.423499999999997
$ws.Range("D69").Value = -7.244099999999992
$ws.Range("D71").Value = -7.595899999999988
$ws.Range("C73").Value = -12.0688
$ws.Range("B76").Value = 6.437599999999994
$ws.Range("A77").Value = -20.00769999999998
$ws.Range("E77").Value = 18.10590000000002
$ws.Range("A78").Value = -20.02679999999998
$ws.Range("A79").Value = -19.88749999999998
$ws.Range("D79").Value = -6.480399999999999
$ws.Range("A80").Value = -19.74119999999998
$ws.Range("A81").Value = -21.6187
$ws.Range("D81").Value = -7.634299999999993
$ws.Range("A82").Value = -22.16890000000001
$ws.Range("E83").Value = 16.899
$ws.Range("A84").Value = -21.99910000000001
$ws.Range("E85").Value = 15.83769999999999
$ws.Range("E86").Value = 16.73280000000001
$ws.Range("B87").Value = 4.179099999999998
$ws.Range("B88").Value = 5.083799999999997
$ws.Range("C90").Value = -12.7035
$ws.Range("B92").Value = 5.589099999999995
$ws.Range("C92").Value = -10.9871
$ws.Range("D92").Value = -6.473099999999996
$ws.Range("C94").Value = -10.05070000000001
$ws.Range("C95").Value = -12.4294
$ws.Range("B96").Value = 5.807099999999999
$ws.Range("E96").Value = 16.40109999999999
$ws.Range("B98").Value = 5.167599999999998
$ws.Range("E98").Value = 16.02889999999999
$ws.Range("A101").Value = -21.20119999999998
$ws.Range("B101").Value = 5.785499999999998
$ws.Range("C101").Value = -12.093
$ws.Range("A102").Value = -21.21919999999999
$ws.Range("B102").Value = 5.4486
